$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6741689443588257
$ws.Range("B1").Value = 2.185449838638306
$ws.Range("C1").Value = 2.381294965744019
$ws.Range("D1").Value = 0.7261258363723755
$ws.Range("E1").Value = 0.8192710280418396
